$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("e2Single")
$ws.Activate()

# Row 45
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = 2
$ws.Range("D45").Value = 3
$ws.Range("F45").Value = 10
$ws.Range("G45").Value = 11
$ws.Range("H45").Value = 12

# Row 46
$ws.Range("B46").Value = 4
$ws.Range("C46").Value = 5
$ws.Range("D46").Value = 6
$ws.Range("F46").Value = 13
$ws.Range("G46").Value = 14
$ws.Range("H46").Value = 15

# Row 48
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = 2
$ws.Range("D48").Value = 3

# Row 49
$ws.Range("B49").Value = 1
$ws.Range("C49").Value = 2
$ws.Range("D49").Value = 3

# Row 50
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 3

# Row 53
$ws.Range("C53").Value = 1
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = 1
$ws.Range("I53").Value = 1

# Row 54
$ws.Range("C54").Value = 2
$ws.Range("F54").Value = 2
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 2

# Row 55
$ws.Range("C55").Value = 3
$ws.Range("E55").Value = 3
$ws.Range("F55").Value = 3
$ws.Range("H55").Value = 3

# Scroll the view down so row 14 is at the top, then select C55 (matches
# the author's final cursor position after entering the new data).
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C55").Select()
